# Update the answer table. Each populated row of the 5-column table gets
# its cell text replaced. Using table-cell addressing (rather than a
# document-wide Find/Replace) avoids collisions between values that are
# simultaneously an "old" value in one cell and a "new" value in another
# (e.g. "28÷5=5, 3" is the old text of row 5 / cell 1, but also becomes
# the new text of row 1 / cell 4).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    # Assigning Range.Text replaces the cell's visible contents while Word
    # keeps the cell's own end-of-cell marker intact.
    $cell.Range.Text = $text
}

# Row 1
Set-CellText $t 1 1 "47÷5=9, 2"
Set-CellText $t 1 2 "66÷7=9, 3"
Set-CellText $t 1 3 "33÷9=3, 6"
Set-CellText $t 1 4 "28÷5=5, 3"
Set-CellText $t 1 5 "82÷6=13, 4"

# Row 5
Set-CellText $t 5 1 "11÷6=1, 5"
Set-CellText $t 5 2 "23÷2=11, 1"
Set-CellText $t 5 3 "38÷7=5, 3"
Set-CellText $t 5 4 "40÷9=4, 4"
Set-CellText $t 5 5 "47÷7=6, 5"

# Row 9
Set-CellText $t 9 1 "76÷2=38, 0"
Set-CellText $t 9 2 "78÷6=13, 0"
Set-CellText $t 9 3 "88÷9=9, 7"
Set-CellText $t 9 4 "75÷5=15, 0"
Set-CellText $t 9 5 "42÷8=5, 2"

# Row 13 (cell 4 keeps its original value; the rest shift/change)
Set-CellText $t 13 1 "94÷6=15, 4"
Set-CellText $t 13 2 "15÷3=5, 0"
Set-CellText $t 13 3 "56÷9=6, 2"
Set-CellText $t 13 4 "59÷3=19, 2"
Set-CellText $t 13 5 "45÷3=15, 0"

# Row 17
Set-CellText $t 17 1 "36÷3=12, 0"
Set-CellText $t 17 2 "86÷9=9, 5"
Set-CellText $t 17 3 "45÷3=15, 0"
Set-CellText $t 17 4 "66÷7=9, 3"
Set-CellText $t 17 5 "72÷5=14, 2"
